$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 15324.214
$ws.Range("I33").Value = 17786.584
$ws.Range("K33").Value = 17786.584
$ws.Range("M33").Value = -17557.584
$ws.Range("H64").Value = 7228.2856
$ws.Range("J64").Value = 6274.5
$ws.Range("L64").Value = 6274.5
$ws.Range("N64").Value = -6770.5
$ws.Range("H67").Value = 7228.2856
$ws.Range("J67").Value = 6274.5
$ws.Range("L67").Value = 6274.5
$ws.Range("N67").Value = -7990.5
$ws.Range("H86").Value = 27782132
$ws.Range("I86").Value = 3707.5
$ws.Range("K86").Value = 3707.5
$ws.Range("M86").Value = -2584.5
$ws.Range("H89").Value = 27782132
$ws.Range("I89").Value = 3707.5
$ws.Range("K89").Value = 18537.5
$ws.Range("M89").Value = -12921.5
$ws.Range("H98").Value = 3286.8
$ws.Range("I98").Value = 3670.6875
$ws.Range("K98").Value = 3670.6875
$ws.Range("M98").Value = -2172.6875
$ws.Range("H122").Value = 3286.8
$ws.Range("I122").Value = 3670.6875
$ws.Range("K122").Value = 11012.0625
$ws.Range("M122").Value = -8562.0625
$ws.Range("H131").Value = 4265.875
$ws.Range("I131").Value = 1906.5555
$ws.Range("J131").Value = 7299.2856
$ws.Range("K131").Value = 5719.666499999999
$ws.Range("L131").Value = 21897.8568
$ws.Range("M131").Value = -679.6664999999994
$ws.Range("N131").Value = -31977.8568
$ws.Range("H132").Value = 1581.037
$ws.Range("I132").Value = 1471.4166
$ws.Range("J132").Value = 2458
$ws.Range("K132").Value = 4414.2498
$ws.Range("L132").Value = 7374
$ws.Range("M132").Value = -1884.2498
$ws.Range("N132").Value = -12434

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H50").Value = 3080.125
$ws.Range("I50").Value = 263
$ws.Range("J50").Value = 4770.4
$ws.Range("K50").Value = 263
$ws.Range("L50").Value = 4770.4
$ws.Range("M50").Value = 451
$ws.Range("N50").Value = -6198.4
$ws.Range("H102").Value = 1616.5714
$ws.Range("I102").Value = 1616.5714
$ws.Range("K102").Value = 1616.5714
$ws.Range("M102").Value = 5.42859999999996

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1673.0454
$ws.Range("I64").Value = 1469.5625
$ws.Range("J64").Value = 2215.6667
$ws.Range("K64").Value = 1469.5625
$ws.Range("L64").Value = 2215.6667
$ws.Range("M64").Value = -1244.5625
$ws.Range("N64").Value = -2665.6667
$ws.Range("H67").Value = 1673.0454
$ws.Range("I67").Value = 1469.5625
$ws.Range("J67").Value = 2215.6667
$ws.Range("K67").Value = 1469.5625
$ws.Range("L67").Value = 2215.6667
$ws.Range("M67").Value = -689.5625
$ws.Range("N67").Value = -3775.6667
$ws.Range("H105").Value = 7353.56
$ws.Range("I105").Value = 8039.952
$ws.Range("J105").Value = 3750
$ws.Range("K105").Value = 8039.952
$ws.Range("L105").Value = 3750
$ws.Range("M105").Value = -6292.952
$ws.Range("N105").Value = -7244

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3677.1177
$ws.Range("I58").Value = 4070
$ws.Range("K58").Value = 4070
$ws.Range("M58").Value = -3867
$ws.Range("H86").Value = 7351.067
$ws.Range("I86").Value = 6643.4287
$ws.Range("J86").Value = 7970.25
$ws.Range("K86").Value = 6643.4287
$ws.Range("L86").Value = 7970.25
$ws.Range("M86").Value = -5520.4287
$ws.Range("N86").Value = -10216.25
$ws.Range("H89").Value = 7351.067
$ws.Range("I89").Value = 6643.4287
$ws.Range("J89").Value = 7970.25
$ws.Range("K89").Value = 33217.14350000001
$ws.Range("L89").Value = 39851.25
$ws.Range("M89").Value = -27601.14350000001
$ws.Range("N89").Value = -51083.25
$ws.Range("H105").Value = 8244.333000000001
$ws.Range("I105").Value = 1449.75
$ws.Range("J105").Value = 21833.5
$ws.Range("K105").Value = 1449.75
$ws.Range("L105").Value = 21833.5
$ws.Range("M105").Value = 297.25
$ws.Range("N105").Value = -25327.5
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("H132").Value = 3041.111
$ws.Range("I132").Value = 2921.375
$ws.Range("K132").Value = 8764.125
$ws.Range("M132").Value = -6234.125
$ws.Range("H134").Value = 4147.5
$ws.Range("I134").Value = 3473.111
$ws.Range("K134").Value = 10419.333
$ws.Range("M134").Value = -7884.332999999999
$ws.Range("H136").Value = 3677.1177
$ws.Range("I136").Value = 4070
$ws.Range("K136").Value = 12210
$ws.Range("M136").Value = -9660
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 666.6667
$ws.Range("I133").Value = 666.6667
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 2000.0001
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = 3059.9999
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6376
$ws.Range("I80").Value = 7752
$ws.Range("J80").Value = 5000
$ws.Range("K80").Value = 7752
$ws.Range("L80").Value = 5000
$ws.Range("M80").Value = -6754
$ws.Range("N80").Value = -6996
$ws.Range("H83").Value = 6376
$ws.Range("I83").Value = 7752
$ws.Range("J83").Value = 5000
$ws.Range("K83").Value = 38760
$ws.Range("L83").Value = 25000
$ws.Range("M83").Value = -33768
$ws.Range("N83").Value = -34984
$ws.Range("H102").Value = 1521.075
$ws.Range("J102").Value = 3399.2222
$ws.Range("L102").Value = 3399.2222
$ws.Range("N102").Value = -6643.2222
$ws.Range("H113").Value = 3301.7058
$ws.Range("J113").Value = 3809.8
$ws.Range("L113").Value = 3809.8
$ws.Range("N113").Value = -8149.8
$ws.Range("H126").Value = 10529725
$ws.Range("I126").Value = 3211.3333
$ws.Range("J126").Value = 20003586
$ws.Range("K126").Value = 9633.999899999999
$ws.Range("L126").Value = 60010758
$ws.Range("M126").Value = -7163.999899999999
$ws.Range("N126").Value = -60015698
$ws.Range("H132").Value = 11760.156
$ws.Range("I132").Value = 11033.12
$ws.Range("J132").Value = 14356.714
$ws.Range("K132").Value = 33099.36
$ws.Range("L132").Value = 43070.142
$ws.Range("M132").Value = -30569.36
$ws.Range("N132").Value = -48130.142
$ws.Range("H136").Value = 82227.86
$ws.Range("J136").Value = 82227.86
$ws.Range("L136").Value = 246683.58
$ws.Range("N136").Value = -251783.58

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 4121.45
$ws.Range("I93").Value = 3772.5293
$ws.Range("K93").Value = 3772.5293
$ws.Range("M93").Value = -2524.5293
$ws.Range("H122").Value = 4281.161
$ws.Range("I122").Value = 3969.9614
$ws.Range("K122").Value = 11909.8842
$ws.Range("M122").Value = -9459.8842

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 6480.4
$ws.Range("I96").Value = 5467.6665
$ws.Range("K96").Value = 5467.6665
$ws.Range("M96").Value = -4094.6665
